$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the bulk of the new journal rows (36-39) first, leaving the two cells
# that get filled in later (D30:D32 rename + D39/E39) for last so the shared
# string table ends up built in the same order as the source workbook.
$ws.Range("A36").Value = "Bruno Díaz"
$ws.Range("B36").Value = 42865
$ws.Range("C36").Value = 2
$ws.Range("E36").Value = "Modificando el MER a partir de la lista de requerimientos"

$ws.Range("A37").Value = "Bruno Díaz"
$ws.Range("B37").Value = 42865
$ws.Range("C37").Value = 3
$ws.Range("E37").Value = "Modificando el MER a partir de la lista de requerimientos (Con ayuda de Fernando)"

$ws.Range("A38").Value = "Bruno Díaz"
$ws.Range("B38").Value = 42868
$ws.Range("C38").Value = 5
$ws.Range("E38").Value = "Posible MER final - Pasaje a MR - Comienzo de creación de Script"

$ws.Range("A39").Value = "Bruno Díaz"
$ws.Range("B39").Value = 42868
$ws.Range("C39").Value = 2
$ws.Range("D39").Value = "Sprint 3 - FrontEnd"

# Rename the task label used for the MER-related entries from
# "Sprint 1 - MER" to "Sprint 3 - Base de Datos" (also applies to the new
# rows above).
$ws.Range("D30").Value = "Sprint 3 - Base de Datos"
$ws.Range("D31").Value = "Sprint 3 - Base de Datos"
$ws.Range("D32").Value = "Sprint 3 - Base de Datos"
$ws.Range("D36").Value = "Sprint 3 - Base de Datos"
$ws.Range("D37").Value = "Sprint 3 - Base de Datos"
$ws.Range("D38").Value = "Sprint 3 - Base de Datos"

$ws.Range("E39").Value = "Copiado de Proyecto de Angular para crear esqueleto - Verificación de que el proyecto corra tal y como estaba para poder realizarle cambios"

# Match the date number format used by the other date cells in column B
# (copy/paste-format instead of assigning a format string so the existing
# style is reused rather than a new numFmt being minted).
$ws.Range("B2").Copy()
$ws.Range("B36:B39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the used range / view to reflect the newly added rows.
$ws.Range("E40").Select()
